# Swap the order of "System" and "dnasr281@gmail.com" within the
# "Recorded By" column (G) for every row that currently reads
# "System, dnasr281@gmail.com", changing it to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
